# Apply updated counts to the "raw" sheet and corresponding percentages to
# the "formatted" sheet of the US issues trends workbook.

$wb = $excel.ActiveWorkbook
$raw = $wb.Worksheets.Item("raw")
$fmt = $wb.Worksheets.Item("formatted")

# ---- "raw" sheet: update raw counts in column C ----
$raw.Range("C3").Value  = 75
$raw.Range("C11").Value = 424
$raw.Range("C14").Value = 225
$raw.Range("C15").Value = 345
$raw.Range("C16").Value = 201
$raw.Range("C19").Value = 221
$raw.Range("C20").Value = 66
$raw.Range("C21").Value = 98
$raw.Range("C24").Value = 61
$raw.Range("C29").Value = 454
$raw.Range("C34").Value = 733
$raw.Range("C35").Value = 34
$raw.Range("C36").Value = 557
$raw.Range("C39").Value = 46
$raw.Range("C47").Value = 88
$raw.Range("C49").Value = 266
$raw.Range("C54").Value = 110
$raw.Range("C55").Value = 555
$raw.Range("C57").Value = 298
$raw.Range("C58").Value = 1353
$raw.Range("C59").Value = 533
$raw.Range("C62").Value = 357
$raw.Range("C63").Value = 110

# ---- "formatted" sheet: update displayed percentage (count) strings ----
$fmt.Range("B2").Value = "39.79% (1,793)"
$fmt.Range("C2").Value = "30.03% (1,353)"
$fmt.Range("D2").Value = "11.83% (533)"

$fmt.Range("B3").Value = "37.79% (1,703)"

$fmt.Range("B4").Value = "16.27% (733)"
$fmt.Range("C4").Value = "16.27% (733)"

$fmt.Range("B5").Value = "12.36% (557)"
$fmt.Range("D5").Value = "12.36% (557)"

$fmt.Range("B6").Value = "12.32% (555)"
$fmt.Range("C6").Value = "12.32% (555)"

$fmt.Range("B7").Value = "10.14% (457)"
$fmt.Range("C7").Value = "7.92% (357)"
$fmt.Range("D7").Value = "2.44% (110)"

$fmt.Range("B8").Value = "10.08% (454)"
$fmt.Range("D8").Value = "10.08% (454)"

$fmt.Range("B9").Value = "9.41% (424)"
$fmt.Range("D9").Value = "9.41% (424)"

$fmt.Range("B10").Value = "7.66% (345)"
$fmt.Range("C10").Value = "7.66% (345)"

$fmt.Range("B11").Value = "6.61% (298)"
$fmt.Range("C11").Value = "6.61% (298)"

$fmt.Range("B12").Value = "5.9% (266)"
$fmt.Range("C12").Value = "5.9% (266)"

$fmt.Range("B13").Value = "4.99% (225)"
$fmt.Range("C13").Value = "4.99% (225)"

$fmt.Range("B14").Value = "4.9% (221)"
$fmt.Range("C14").Value = "4.9% (221)"

$fmt.Range("B15").Value = "4.46% (201)"
$fmt.Range("C15").Value = "4.46% (201)"

$fmt.Range("B16").Value = "1.18% (53)"
